$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at K (before former "short_name") for "categories"
$ws.Columns.Item(11).Insert()

# Insert two new columns at X:Y (before former "codelist") for
# "derived_variable" and "derivation_description"
$ws.Range("X1:Y1").EntireColumn.Insert()

# Populate the three new header cells (shared strings are appended
# automatically / style is inherited from the shift, matching s="2" for K1
# and s="1" for X1/Y1)
$ws.Range("K1").Value = "categories"
$ws.Range("X1").Value = "derived_variable"
$ws.Range("Y1").Value = "derivation_description"

# Restore explicit column widths for the newly inserted columns
$ws.Columns.Item(11).ColumnWidth = 21.142857142857142
$ws.Columns.Item(24).ColumnWidth = 18.714285714285715
$ws.Columns.Item(25).ColumnWidth = 38.714285714285715

# Rebuild the AutoFilter over the new full range A1:AI1
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:AI1").AutoFilter()

# Update the workbook-level _FilterDatabase defined name to match
$name = $wb.Names.Item(1)
$name.RefersTo = "='Collection Specializations'!`$A`$1:`$AI`$1"

# Reflect the view state recorded in the saved file (active cell moved to Y2)
[void]$ws.Range("Y2").Select()
